$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain plain text (matching source inlineStr cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.816.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.813.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.43"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000251"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.13"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.450.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.822.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.837.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.84"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000148"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.12"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.963.09"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.43"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.48"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.25"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.21"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.300"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.56"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.85"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "392.15"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.31%  "
